# "TestNG And Advanced Actions"
#
# The author duplicated the "Product_Info" worksheet (creating a second,
# identical copy named "Product_Info2") and placed the new copy directly
# in front of the original "Product_Info" tab, giving a final tab order of:
#   Login-Info, Product_Info2, Product_Info
# The new sheet kept the same data/shared-string values as "Product_Info".
# The author was left with the new "Product_Info2" sheet active/selected,
# while the other two sheets' selections also moved around a bit.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate "Product_Info", inserting the copy right before it ---
$srcSheet = $wb.Worksheets.Item("Product_Info")
$srcSheet.Copy($srcSheet, $null)

# Excel names the copy "Product_Info (2)" and drops it immediately before
# the original -> it is now worksheet #2 ("Login-Info" is #1). Rename it.
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "Product_Info2"

$wsLogin    = $wb.Worksheets.Item("Login-Info")
$wsProduct2 = $wb.Worksheets.Item("Product_Info2")
$wsProduct  = $wb.Worksheets.Item("Product_Info")

# --- 2. Column-width tweaks made on the new "Product_Info2" copy ---
# (only the first three columns were resized - the fourth was left as-is)
$wsProduct2.Columns.Item(1).ColumnWidth = 19
$wsProduct2.Columns.Item(2).ColumnWidth = 21.333333333333332
$wsProduct2.Columns.Item(3).ColumnWidth = 25.333333333333332

# --- 3. Selections left behind on each sheet when the file was saved ---
[void]$wsLogin.Range("I5").Select()
[void]$wsProduct.Range("B1:D5").Select()

# "Product_Info2" is the sheet that was active/selected when saved, so
# activate it last and select its last-used cell.
$wsProduct2.Activate()
[void]$wsProduct2.Range("B22").Select()
